$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> number map (A..H)
$colNum = @{ 'A' = 1; 'B' = 2; 'C' = 3; 'D' = 4; 'E' = 5; 'F' = 6; 'G' = 7; 'H' = 8 }

# Row-level updates: each entry maps a row number to the column/value pairs
# that changed for that row in this data refresh (new countries inserted,
# updated case counts, and the "last updated" timestamp).
$updates = @{
    1   = @{ A = "Datos actualizados a 12 de Abril de 2020 a las 01:52" }
    4   = @{ B = 531943; C = 29067; D = 29818; E = 481570 }
    8   = @{ B = 125452; C = 3281; E = 65181; G = 135; H = 2871 }
    17  = @{ B = 20962; C = 1173; E = 19649; G = 72; H = 1140 }
    23  = @{ F = 789 }
    45  = @{ B = 3234; C = 260; D = 23; E = 3132; G = 5; H = 79 }
    144 = @{ A = "Islas Caimanes"; B = 53; C = 8; D = 6; E = 46; F = 3; H = 1 }
    145 = @{ A = "Polinesia Francesa"; B = 51; D = 0; E = 51; F = 0; G = 0; H = 0 }
    146 = @{ A = "San Martin (Parte Holandesa)"; B = 50; C = 0; D = 5; E = 36; F = 2; G = 1; H = 9 }
    147 = @{ A = "Liberia"; C = 11; D = 3; E = 40; F = 0; H = 5 }
    148 = @{ A = "Bermudas"; B = 48; C = 0; D = 25; E = 19; F = 2; H = 4 }
    149 = @{ A = "Gabon"; C = 2; D = 1; E = 44; F = 0; H = 1 }
    150 = @{ A = "Bahamas"; B = 46; C = 4; D = 5; E = 33; F = 1; H = 8 }
    164 = @{ B = 25; C = 1; E = 16 }
    194 = @{ A = "Islas Turcas y Caicos"; C = 1; E = 8; H = 1 }
    195 = @{ A = "Nicaragua"; C = 2; D = 0; E = 8 }
    196 = @{ A = "Montserrat"; B = 9; E = 7; H = 2 }
    197 = @{ A = "Gambia"; B = 9; C = 5; D = 2; E = 6; H = 1 }
    198 = @{ A = "Sierra Leona"; E = 8; H = 0 }
    199 = @{ A = "Republica de Africa Central"; C = 0; E = 8; H = 0 }
    200 = @{ A = "Santa Sede"; C = 0; D = 2; H = 0 }
    201 = @{ A = "Cabo Verde"; C = 1; D = 1; H = 1 }
    203 = @{ A = "Sahara Occidental"; C = 2; D = 0; E = 6 }
    204 = @{ A = "San Bartolome"; B = 6; C = 0; D = 1 }
    205 = @{ A = "Burundi"; C = 2; D = 0; E = 5 }
    206 = @{ A = "Islas Malvinas"; D = 1; E = 4 }
    207 = @{ A = "Butan"; B = 5; D = 2; E = 3 }
    208 = @{ A = "Santo Tome y Principe" }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item($rowNum, $colNum[$col]).Value = $cols[$col]
    }
}
